$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-08-11 Monday" "2025-08-12 Tuesday"

Replace-Text "425÷4=106, 1" "452÷4=113, 0"
Replace-Text "843÷4=210, 3" "627÷2=313, 1"
Replace-Text "520÷8=65, 0" "433÷2=216, 1"
Replace-Text "224÷2=112, 0" "236÷5=47, 1"
Replace-Text "937÷3=312, 1" "458÷3=152, 2"

Replace-Text "895÷6=149, 1" "573÷9=63, 6"
Replace-Text "533÷4=133, 1" "797÷2=398, 1"
Replace-Text "383÷5=76, 3" "586÷6=97, 4"
Replace-Text "238÷3=79, 1" "707÷8=88, 3"
Replace-Text "604÷3=201, 1" "478÷6=79, 4"

Replace-Text "579÷8=72, 3" "210÷2=105, 0"
Replace-Text "110÷4=27, 2" "429÷3=143, 0"
Replace-Text "413÷6=68, 5" "232÷5=46, 2"
Replace-Text "935÷8=116, 7" "201÷9=22, 3"
Replace-Text "342÷4=85, 2" "795÷4=198, 3"

Replace-Text "220÷4=55, 0" "280÷7=40, 0"
Replace-Text "453÷7=64, 5" "961÷6=160, 1"
Replace-Text "146÷2=73, 0" "723÷6=120, 3"
Replace-Text "345÷6=57, 3" "795÷7=113, 4"
Replace-Text "538÷7=76, 6" "332÷9=36, 8"

Replace-Text "114÷4=28, 2" "771÷7=110, 1"
Replace-Text "793÷4=198, 1" "152÷3=50, 2"
Replace-Text "938÷8=117, 2" "789÷6=131, 3"
Replace-Text "535÷4=133, 3" "279÷4=69, 3"
Replace-Text "677÷3=225, 2" "625÷7=89, 2"
